# Daily attendance processing - 2025-12-05 22:50:28
# For every row in the "Recorded By" column (G), when the value is a
# two-part, comma-separated list with "System" listed first
# (e.g. "System, dnasr281@gmail.com"), swap the order so "System" is
# listed last (e.g. "dnasr281@gmail.com, System"). Values with a
# different shape (single name, 3+ names, or "System" not first) are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Text

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ', '
        if ($parts.Count -eq 2 -and $parts[0] -eq 'System') {
            $cell.Value = $parts[1] + ", " + $parts[0]
        }
    }
}
